$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Formatting -----------------------------------------------------

# Column A ("id") reuses the bold / thin-bordered / centered look already
# used for the rest of the "id" column. Build it once on the first new
# row, then fan it out to the remaining new rows.
$idFirst = $ws.Range("A165")
$idFirst.HorizontalAlignment = -4108
$idFirst.VerticalAlignment = -4160
$idFirst.Borders.LineStyle = 1
$idFirst.Font.Bold = $true
$idFirst.Copy()
$ws.Range("A166:A172").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Scratch cell (outside the table, cleared at the end): format as Text so
# values like "12/03/2018" or "0"/"1" round-trip as literal strings, then
# paste-special VALUES ONLY into the real destination cell so the
# destination keeps the table's normal (General) formatting, matching the
# rest of the sheet where the text cells carry no explicit number format.
$scratch = $ws.Range("AZ1")
$scratch.NumberFormat = "@"

# Row 165 (id 164)
$ws.Range("A165").Value = 164
$ws.Range("B165").Value = "https://github.com/tensorflow/ranking"
$ws.Range("C165").Value = "ranking"
$ws.Range("D165").Value = "tensorflow"
$scratch.Value = "12/03/2018"
$scratch.Copy()
$ws.Range("E165").PasteSpecial(-4163)
$scratch.Value = "0"
$scratch.Copy()
$ws.Range("F165").PasteSpecial(-4163)
$scratch.Value = "0"
$scratch.Copy()
$ws.Range("G165").PasteSpecial(-4163)
$scratch.Value = "0"
$scratch.Copy()
$ws.Range("H165").PasteSpecial(-4163)
$scratch.Value = "1"
$scratch.Copy()
$ws.Range("I165").PasteSpecial(-4163)
$scratch.Value = "1"
$scratch.Copy()
$ws.Range("J165").PasteSpecial(-4163)
$scratch.Value = "1"
$scratch.Copy()
$ws.Range("K165").PasteSpecial(-4163)
$scratch.Value = "0"
$scratch.Copy()
$ws.Range("L165").PasteSpecial(-4163)
$scratch.Value = "0"
$scratch.Copy()
$ws.Range("M165").PasteSpecial(-4163)
$scratch.Value = "1"
$scratch.Copy()
$ws.Range("N165").PasteSpecial(-4163)
$scratch.Value = "0"
$scratch.Copy()
$ws.Range("O165").PasteSpecial(-4163)

# Row 166 (id 165)
$ws.Range("A166").Value = 165
$ws.Range("B166").Value = "https://github.com/tensorflow/ranking"
$ws.Range("C166").Value = "ranking"
$ws.Range("D166").Value = "tensorflow"
$scratch.Value = "01/12/2020"
$scratch.Copy()
$ws.Range("E166").PasteSpecial(-4163)
$scratch.Value = "0"
$scratch.Copy()
$ws.Range("F166").PasteSpecial(-4163)
$scratch.Value = "0"
$scratch.Copy()
$ws.Range("G166").PasteSpecial(-4163)
$scratch.Value = "0"
$scratch.Copy()
$ws.Range("H166").PasteSpecial(-4163)
$scratch.Value = "1"
$scratch.Copy()
$ws.Range("I166").PasteSpecial(-4163)
$scratch.Value = "0"
$scratch.Copy()
$ws.Range("J166").PasteSpecial(-4163)
$scratch.Value = "0"
$scratch.Copy()
$ws.Range("K166").PasteSpecial(-4163)
$scratch.Value = "0"
$scratch.Copy()
$ws.Range("L166").PasteSpecial(-4163)
$scratch.Value = "0"
$scratch.Copy()
$ws.Range("M166").PasteSpecial(-4163)
$scratch.Value = "1"
$scratch.Copy()
$ws.Range("N166").PasteSpecial(-4163)
$scratch.Value = "0"
$scratch.Copy()
$ws.Range("O166").PasteSpecial(-4163)

# Row 167 (id 166)
$ws.Range("A167").Value = 166
$ws.Range("B167").Value = "https://github.com/tensorflow/ranking"
$ws.Range("C167").Value = "ranking"
$ws.Range("D167").Value = "tensorflow"
$scratch.Value = "12/03/2018"
$scratch.Copy()
$ws.Range("E167").PasteSpecial(-4163)
$scratch.Value = "0"
$scratch.Copy()
$ws.Range("F167").PasteSpecial(-4163)
$scratch.Value = "0"
$scratch.Copy()
$ws.Range("G167").PasteSpecial(-4163)
$scratch.Value = "0"
$scratch.Copy()
$ws.Range("H167").PasteSpecial(-4163)
$scratch.Value = "1"
$scratch.Copy()
$ws.Range("I167").PasteSpecial(-4163)
$scratch.Value = "1"
$scratch.Copy()
$ws.Range("J167").PasteSpecial(-4163)
$scratch.Value = "1"
$scratch.Copy()
$ws.Range("K167").PasteSpecial(-4163)
$scratch.Value = "0"
$scratch.Copy()
$ws.Range("L167").PasteSpecial(-4163)
$scratch.Value = "0"
$scratch.Copy()
$ws.Range("M167").PasteSpecial(-4163)
$scratch.Value = "1"
$scratch.Copy()
$ws.Range("N167").PasteSpecial(-4163)
$scratch.Value = "0"
$scratch.Copy()
$ws.Range("O167").PasteSpecial(-4163)

# Row 168 (id 167)
$ws.Range("A168").Value = 167
$ws.Range("B168").Value = "https://github.com/tensorflow/ranking"
$ws.Range("C168").Value = "ranking"
$ws.Range("D168").Value = "tensorflow"
$scratch.Value = "12/03/2018"
$scratch.Copy()
$ws.Range("E168").PasteSpecial(-4163)
$scratch.Value = "0"
$scratch.Copy()
$ws.Range("F168").PasteSpecial(-4163)
$scratch.Value = "0"
$scratch.Copy()
$ws.Range("G168").PasteSpecial(-4163)
$scratch.Value = "0"
$scratch.Copy()
$ws.Range("H168").PasteSpecial(-4163)
$scratch.Value = "1"
$scratch.Copy()
$ws.Range("I168").PasteSpecial(-4163)
$scratch.Value = "1"
$scratch.Copy()
$ws.Range("J168").PasteSpecial(-4163)
$scratch.Value = "1"
$scratch.Copy()
$ws.Range("K168").PasteSpecial(-4163)
$scratch.Value = "0"
$scratch.Copy()
$ws.Range("L168").PasteSpecial(-4163)
$scratch.Value = "0"
$scratch.Copy()
$ws.Range("M168").PasteSpecial(-4163)
$scratch.Value = "1"
$scratch.Copy()
$ws.Range("N168").PasteSpecial(-4163)
$scratch.Value = "0"
$scratch.Copy()
$ws.Range("O168").PasteSpecial(-4163)

# Row 169 (id 168)
$ws.Range("A169").Value = 168
$ws.Range("B169").Value = "https://github.com/tensorflow/ranking"
$ws.Range("C169").Value = "ranking"
$ws.Range("D169").Value = "tensorflow"
$scratch.Value = "01/12/2020"
$scratch.Copy()
$ws.Range("E169").PasteSpecial(-4163)
$scratch.Value = "0"
$scratch.Copy()
$ws.Range("F169").PasteSpecial(-4163)
$scratch.Value = "0"
$scratch.Copy()
$ws.Range("G169").PasteSpecial(-4163)
$scratch.Value = "0"
$scratch.Copy()
$ws.Range("H169").PasteSpecial(-4163)
$scratch.Value = "1"
$scratch.Copy()
$ws.Range("I169").PasteSpecial(-4163)
$scratch.Value = "0"
$scratch.Copy()
$ws.Range("J169").PasteSpecial(-4163)
$scratch.Value = "0"
$scratch.Copy()
$ws.Range("K169").PasteSpecial(-4163)
$scratch.Value = "0"
$scratch.Copy()
$ws.Range("L169").PasteSpecial(-4163)
$scratch.Value = "0"
$scratch.Copy()
$ws.Range("M169").PasteSpecial(-4163)
$scratch.Value = "1"
$scratch.Copy()
$ws.Range("N169").PasteSpecial(-4163)
$scratch.Value = "0"
$scratch.Copy()
$ws.Range("O169").PasteSpecial(-4163)

# Row 170 (id 169)
$ws.Range("A170").Value = 169
$ws.Range("B170").Value = "https://github.com/tensorflow/ranking"
$ws.Range("C170").Value = "ranking"
$ws.Range("D170").Value = "tensorflow"
$scratch.Value = "01/12/2020"
$scratch.Copy()
$ws.Range("E170").PasteSpecial(-4163)
$scratch.Value = "0"
$scratch.Copy()
$ws.Range("F170").PasteSpecial(-4163)
$scratch.Value = "0"
$scratch.Copy()
$ws.Range("G170").PasteSpecial(-4163)
$scratch.Value = "0"
$scratch.Copy()
$ws.Range("H170").PasteSpecial(-4163)
$scratch.Value = "1"
$scratch.Copy()
$ws.Range("I170").PasteSpecial(-4163)
$scratch.Value = "0"
$scratch.Copy()
$ws.Range("J170").PasteSpecial(-4163)
$scratch.Value = "0"
$scratch.Copy()
$ws.Range("K170").PasteSpecial(-4163)
$scratch.Value = "0"
$scratch.Copy()
$ws.Range("L170").PasteSpecial(-4163)
$scratch.Value = "0"
$scratch.Copy()
$ws.Range("M170").PasteSpecial(-4163)
$scratch.Value = "1"
$scratch.Copy()
$ws.Range("N170").PasteSpecial(-4163)
$scratch.Value = "0"
$scratch.Copy()
$ws.Range("O170").PasteSpecial(-4163)

# Row 171 (id 170)
$ws.Range("A171").Value = 170
$ws.Range("B171").Value = "https://github.com/tensorflow/ranking"
$ws.Range("C171").Value = "ranking"
$ws.Range("D171").Value = "tensorflow"
$scratch.Value = "01/12/2020"
$scratch.Copy()
$ws.Range("E171").PasteSpecial(-4163)
$scratch.Value = "0"
$scratch.Copy()
$ws.Range("F171").PasteSpecial(-4163)
$scratch.Value = "0"
$scratch.Copy()
$ws.Range("G171").PasteSpecial(-4163)
$scratch.Value = "0"
$scratch.Copy()
$ws.Range("H171").PasteSpecial(-4163)
$scratch.Value = "1"
$scratch.Copy()
$ws.Range("I171").PasteSpecial(-4163)
$scratch.Value = "0"
$scratch.Copy()
$ws.Range("J171").PasteSpecial(-4163)
$scratch.Value = "0"
$scratch.Copy()
$ws.Range("K171").PasteSpecial(-4163)
$scratch.Value = "0"
$scratch.Copy()
$ws.Range("L171").PasteSpecial(-4163)
$scratch.Value = "0"
$scratch.Copy()
$ws.Range("M171").PasteSpecial(-4163)
$scratch.Value = "1"
$scratch.Copy()
$ws.Range("N171").PasteSpecial(-4163)
$scratch.Value = "0"
$scratch.Copy()
$ws.Range("O171").PasteSpecial(-4163)

# Row 172 (id 171)
$ws.Range("A172").Value = 171
$ws.Range("B172").Value = "https://github.com/tensorflow/ranking"
$ws.Range("C172").Value = "ranking"
$ws.Range("D172").Value = "tensorflow"
$scratch.Value = "12/03/2018"
$scratch.Copy()
$ws.Range("E172").PasteSpecial(-4163)
$scratch.Value = "0"
$scratch.Copy()
$ws.Range("F172").PasteSpecial(-4163)
$scratch.Value = "0"
$scratch.Copy()
$ws.Range("G172").PasteSpecial(-4163)
$scratch.Value = "0"
$scratch.Copy()
$ws.Range("H172").PasteSpecial(-4163)
$scratch.Value = "1"
$scratch.Copy()
$ws.Range("I172").PasteSpecial(-4163)
$scratch.Value = "1"
$scratch.Copy()
$ws.Range("J172").PasteSpecial(-4163)
$scratch.Value = "1"
$scratch.Copy()
$ws.Range("K172").PasteSpecial(-4163)
$scratch.Value = "0"
$scratch.Copy()
$ws.Range("L172").PasteSpecial(-4163)
$scratch.Value = "0"
$scratch.Copy()
$ws.Range("M172").PasteSpecial(-4163)
$scratch.Value = "1"
$scratch.Copy()
$ws.Range("N172").PasteSpecial(-4163)
$scratch.Value = "0"
$scratch.Copy()
$ws.Range("O172").PasteSpecial(-4163)

# Remove the scratch cell so it does not show up in the used range.
$scratch.Clear()
$excel.CutCopyMode = $false